# Add a new row to the "Users" sheet for a user used to test annotation
# sharing with a group (mirrors the existing annotation-user rows 19-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Copy the formatting of the row above (border etc.) down into row 22 first,
# then fill in the values - this matches how the row was authored in Excel.
$ws.Range("A19:G19").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)

# G22 additionally carries the "Hyperlink" look used by the two rows above it.
$ws.Range("G21").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("A22").Value = "userForAnnotationGroup"
$ws.Range("B22").Value = "Password1"
$ws.Range("E22").Value = "ANZ annotation user"
$ws.Range("G22").Value = "userforannotationgroup@mailinator.com"

$ws.Range("E22").Select()
